$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = 'Normal'
}

Set-TextCell 'D2' '28.483.16'
Set-TextCell 'E2' '  +0.07%  '
Set-TextCell 'D3' '1.826.17'
Set-TextCell 'E3' '  -0.05%  '
Set-TextCell 'E4' '  +0.25%  '
Set-TextCell 'D5' '316.36'
Set-TextCell 'E5' '  +0.30%  '
Set-TextCell 'E6' '  +0.18%  '
Set-TextCell 'D7' '0.5161'
Set-TextCell 'E7' '  +2.05%  '
Set-TextCell 'D8' '0.3864'
Set-TextCell 'E8' '  -1.21%  '
Set-TextCell 'D9' '0.08291'
Set-TextCell 'E9' '  +8.23%  '
Set-TextCell 'D11' '41.93'
Set-TextCell 'E11' '  +0.06%  '
Set-TextCell 'D12' '6.389'
Set-TextCell 'E12' '  +1.63%  '
Set-TextCell 'D13' '21.20'
Set-TextCell 'E13' '  +0.63%  '
Set-TextCell 'E14' '  +0.20%  '
Set-TextCell 'D15' '7.492'
Set-TextCell 'E15' '  -0.95%  '
Set-TextCell 'D16' '1.828.78'
Set-TextCell 'E16' '  +0.20%  '
Set-TextCell 'D17' '94.03'
Set-TextCell 'E17' '  +0.89%  '
Set-TextCell 'D18' '0.00001123'
Set-TextCell 'E18' '  +3.59%  '
Set-TextCell 'E19' '  -0.32%  '
Set-TextCell 'D20' '17.80'
Set-TextCell 'E20' '  +0.59%  '
Set-TextCell 'E21' '  +0.19%  '
Set-TextCell 'D22' '6.062'
Set-TextCell 'E22' '  -1.60%  '
Set-TextCell 'D23' '28.521.44'
Set-TextCell 'E23' '  +0.09%  '
Set-TextCell 'E24' '  +2.84%  '
Set-TextCell 'D25' '2.244'
Set-TextCell 'E25' '  -0.56%  '
Set-TextCell 'D26' '21.09'
Set-TextCell 'E26' '  +2.34%  '
Set-TextCell 'D27' '159.23'
Set-TextCell 'E27' '  +1.62%  '
Set-TextCell 'D28' '2.038.21'
Set-TextCell 'E28' '  +0.05%  '
Set-TextCell 'D29' '2.419'
Set-TextCell 'E29' '  +1.06%  '
Set-TextCell 'D30' '125.85'
Set-TextCell 'E30' '  +0.62%  '
Set-TextCell 'D31' '0.1095'
Set-TextCell 'E31' '  +1.29%  '
Set-TextCell 'D32' '1.098'
Set-TextCell 'E32' '  -2.65%  '
Set-TextCell 'D33' '0.07662'
Set-TextCell 'E33' '  +8.61%  '
Set-TextCell 'D34' '5.730'
Set-TextCell 'E34' '  +0.98%  '
Set-TextCell 'D35' '3.681'
Set-TextCell 'E35' '  +0.51%  '
Set-TextCell 'D36' '0.2238'
Set-TextCell 'E36' '  +0.70%  '
Set-TextCell 'D37' '0.02372'
Set-TextCell 'E37' '  +2.18%  '
Set-TextCell 'D38' '5.269'
Set-TextCell 'E38' '  +2.55%  '
Set-TextCell 'D39' '12.02'
Set-TextCell 'E39' '  +7.11%  '
Set-TextCell 'D40' '8.781'
Set-TextCell 'E40' '  -1.92%  '
Set-TextCell 'D41' '0.6434'
Set-TextCell 'E41' '  +3.01%  '
Set-TextCell 'D42' '1.192'
Set-TextCell 'E42' '  +0.83%  '
Set-TextCell 'D43' '1.401'
Set-TextCell 'E43' '  +0.21%  '
Set-TextCell 'B44' 'Decentraland'
Set-TextCell 'C44' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 'D44' '0.6225'
Set-TextCell 'E44' '  +5.53%  '
Set-TextCell 'B45' 'EnergySwap'
Set-TextCell 'C45' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D45' '13.59'
Set-TextCell 'E45' '  +1.43%  '
Set-TextCell 'D46' '3.800'
Set-TextCell 'E46' '  +2.19%  '
Set-TextCell 'D47' '127.90'
Set-TextCell 'E47' '  +2.77%  '
Set-TextCell 'E48' '  +1.17%  '
Set-TextCell 'D49' '1.204'
Set-TextCell 'E49' '  +0.85%  '
Set-TextCell 'D50' '0.06976'
Set-TextCell 'E50' '  +0.93%  '
Set-TextCell 'D51' '74.27'
Set-TextCell 'E51' '  +0.53%  '
